# Apply the edits described by the commit "1026: Fixed failed unit tests":
#  - Rename the worksheet from "Basic Clinic Data..." to
#    "Test Import Survey Response 1" (this is the sheet's test-fixture name
#    used for importing a single survey's responses).
#  - Move/replace the saved cell selection on the sheet from G12 to F20.
#
# (Note: the diff also shows changes to the workbook's internal
#  x15ac:absPath and xr:revisionPtr/documentId metadata. Those are
#  Microsoft Office authoring/telemetry artifacts written automatically by
#  the real Excel application based on the machine/session that saved the
#  file - they are not exposed anywhere in the Excel COM object model and
#  cannot be set from script, so they are intentionally left untouched
#  here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet.
$ws.Name = "Test Import Survey Response 1"

# Update the active selection saved in the sheet view.
$ws.Range("F20").Select()
